$d = $word.ActiveDocument
$wordMlNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# 1. After "...Svakom korisniku će biti omogućeno da kreira, uređuje i briše
#    svoj događaj unutar aplikacije." append a new sentence/run to the same
#    paragraph, then add a brand-new paragraph with a second sentence, right
#    before the existing empty paragraph that follows.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Svakom korisniku će biti omogućeno da kreira, uređuje i briše svoj događaj unutar aplikacije.")
$rng.Collapse(0)
$insertStart = $rng.Start
$sentence1 = " Za svaki događaj koji korisnik kreira će morati da unese potrebne podatke i da podesi njegovu vidljivost. "
$rng.InsertAfter($sentence1)

# Force the freshly inserted text to live in its own <w:r> (rather than being
# merged into the preceding run) by toggling a character formatting property
# off/on across the whole inserted span.
$newRunRng = $d.Range($insertStart, $insertStart + $sentence1.Length)
$newRunRng.Bold = 1
$newRunRng.Bold = 0

# Insert a paragraph break right after the text we just added, then type the
# second sentence into the newly created paragraph.
$afterIns = $insertStart + $sentence1.Length
$breakPoint = $d.Range($afterIns, $afterIns)
$breakPoint.InsertParagraphAfter()

$sentence2 = "Svaki događaj će imati svoju vidljivost (podrazumevano je javno) koja određuje domen korisnika koji mogu da ga vide i da se prijave na njega."
$para2Start = $afterIns + 1
$para2Rng = $d.Range($para2Start, $para2Start)
$para2Rng.InsertAfter($sentence2)

# ---------------------------------------------------------------------------
# 2. Drop the stray "_GoBack" bookmark and let the two runs it used to sit
#    between ("<space>" and "kao i lakše upravljanje...") collapse into one
#    run, while the runs on either side (different rsid) stay untouched.
#    Trick: temporarily relocate the _GoBack bookmark (and add a scratch
#    bookmark) so they act as barriers that keep the outer runs from being
#    swept into the merge; then do a no-op edit inside the target span to
#    trigger the run-coalescing pass; finally remove both bookmarks.
# ---------------------------------------------------------------------------
$rngLeft = $d.Content
$rngLeft.Find.Execute("preko olakšane distribucije pozivnica")
$leftBoundary = $rngLeft.End

$rngRight = $d.Content
$rngRight.Find.Execute("olakšana komunikacija sa njima.")
$rightBoundary = $rngRight.Start

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$rightAnchor = $d.Range($rightBoundary, $rightBoundary)
$d.Bookmarks.Add("_GoBack", $rightAnchor)

$leftAnchor = $d.Range($leftBoundary, $leftBoundary)
$d.Bookmarks.Add("ZZScratchBarrier", $leftAnchor)

$rngMid = $d.Content
$rngMid.Find.Execute("kao i lakše")
$tinyPoint = $d.Range($rngMid.Start + 1, $rngMid.Start + 1)
$tinyPoint.InsertAfter("x")
$tinyDel = $d.Range($rngMid.Start + 1, $rngMid.Start + 2)
$tinyDel.Delete()

if ($d.Bookmarks.Exists("ZZScratchBarrier")) {
    $d.Bookmarks("ZZScratchBarrier").Delete()
}
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3/4/5. Relocate the lastRenderedPageBreak marker: remove it from
#    "Tehnologije u upotrebi" and from "Za vođu tima se bira...", and add it
#    to the "Adnan Crnovršanin" run that follows the "Vođa tima:" label.
# ---------------------------------------------------------------------------
$rngTech = $d.Content
$rngTech.Find.Execute("Tehnologije u upotrebi")
$techStart = $rngTech.Start
$rngTech.Delete()
$techCollapsed = $d.Range($techStart, $techStart)
$techCollapsed.InsertXML("<w:p xmlns:w='$wordMlNs'><w:r><w:rPr><w:lang w:val='sr-Latn-RS'/></w:rPr><w:t>Tehnologije u upotrebi</w:t></w:r></w:p>")

$rngVoda = $d.Content
$rngVoda.Find.Execute("Vođa tima:")
$rngAdnan = $d.Range($rngVoda.End, $d.Content.End)
$rngAdnan.Find.Execute("Adnan Crnovršanin")
$adnanStart = $rngAdnan.Start
$rngAdnan.Delete()
$adnanCollapsed = $d.Range($adnanStart, $adnanStart)
$adnanCollapsed.InsertXML("<w:p xmlns:w='$wordMlNs'><w:r><w:rPr><w:lang w:val='sr-Latn-RS'/></w:rPr><w:lastRenderedPageBreak/><w:t>Adnan Crnovršanin</w:t></w:r></w:p>")

$rngZa = $d.Content
$rngZa.Find.Execute("Za vođu tima se bira Adnan Crnovršanin radi sticanja iskustva i učenja odgovornosti i obaveza ove pozicije u ovakvom jednom kontrolisanom okruženju.")
$zaStart = $rngZa.Start
$rngZa.Delete()
$zaCollapsed = $d.Range($zaStart, $zaStart)
$zaCollapsed.InsertXML("<w:p xmlns:w='$wordMlNs'><w:r><w:rPr><w:lang w:val='sr-Latn-RS'/></w:rPr><w:t>Za vođu tima se bira Adnan Crnovršanin radi sticanja iskustva i učenja odgovornosti i obaveza ove pozicije u ovakvom jednom kontrolisanom okruženju.</w:t></w:r></w:p>")

Write-Output "All edits applied"
